$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.325.09'
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('D3').Value = '2.098.17'
$ws.Range('E3').Value = '  +3.07%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '228.75'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  +0.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.30'
$ws.Range('E7').Value = '  +1.67%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.380'
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0847'
$ws.Range('E10').Value = '  +3.01%  '
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').Value = '2.413.61'
$ws.Range('E12').Value = '  +3.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.74'
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.35'
$ws.Range('E14').Value = '  +6.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.48'
$ws.Range('E15').Value = '  +5.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.777'
$ws.Range('E16').Value = '  +2.21%  '
$ws.Range('D17').Value = '2.120.53'
$ws.Range('E17').Value = '  +2.47%  '
$ws.Range('D18').Value = '38.263.25'
$ws.Range('E18').Value = '  +1.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.01'
$ws.Range('E19').Value = '  +1.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.39'
$ws.Range('E20').Value = '  +0.84%  '
$ws.Range('D21').Value = '0.0₃0833'
$ws.Range('E21').Value = '  +1.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '224.63'
$ws.Range('E22').Value = '  +0.37%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('E25').Value = '  +2.50%  '
$ws.Range('E26').Value = '  +1.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.42'
$ws.Range('E27').Value = '  +0.42%  '
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.01'
$ws.Range('E29').Value = '  +0.75%  '
$ws.Range('E30').Value = '  +6.46%  '
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('E32').Value = '  +9.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.72'
$ws.Range('E33').Value = '  +4.87%  '
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0605'
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('E36').Value = '  +4.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.40'
$ws.Range('E37').Value = '  +1.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.54'
$ws.Range('E38').Value = '  +6.30%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '18.14'
$ws.Range('E40').Value = '  +2.23%  '
$ws.Range('D41').Value = '1.545.93'
$ws.Range('E41').Value = '  +0.60%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.90'
$ws.Range('E42').Value = '  +3.84%  '
$ws.Range('E43').Value = '  +0.99%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.83'
$ws.Range('E44').Value = '  +1.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0909'
$ws.Range('E45').Value = '  -0.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.18'
$ws.Range('E46').Value = '  +3.97%  '
$ws.Range('E47').Value = '  +1.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.52'
$ws.Range('E48').Value = '  +5.05%  '
$ws.Range('E49').Value = '  +1.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.99'
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('D51').Value = '2.299.62'
